$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F7").Value = 1019
$ws.Range("F8").Value = 1669
$ws.Range("F12").Value = 370
$ws.Range("F15").Value = 1185
$ws.Range("F18").Value = 2123
$ws.Range("F22").Value = 565
$ws.Range("F23").Value = 12
$ws.Range("F24").Value = 1353
$ws.Range("F25").Value = 1131
$ws.Range("F33").Value = 1180
$ws.Range("F36").Value = 911
$ws.Range("F41").Value = 137
$ws.Range("F42").Value = 2127
$ws.Range("F45").Value = 1176

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F11").Value = 284
$ws.Range("F26").Value = 201

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 2970
$ws.Range("F10").Value = 837
$ws.Range("F12").Value = 497
$ws.Range("F13").Value = 1219
$ws.Range("F15").Value = 894

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F7").Value = 837
$ws.Range("F10").Value = 497
$ws.Range("F11").Value = 1219
$ws.Range("F12").Value = 1019
$ws.Range("F13").Value = 1669
$ws.Range("F17").Value = 284
$ws.Range("F20").Value = 1185
$ws.Range("F21").Value = 894
$ws.Range("F22").Value = 894
$ws.Range("F23").Value = 2123
$ws.Range("F28").Value = 565
$ws.Range("F29").Value = 1353
$ws.Range("F31").Value = 1131
$ws.Range("F39").Value = 1180
$ws.Range("F41").Value = 911
$ws.Range("F46").Value = 2127
$ws.Range("F49").Value = 1176
